$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header updates: bulletin number and report week dates ---
$ws.Range("A8").Value = "Volume 31   Number  43"
$ws.Range("C9").Value = "Report Covering the Week  10/21/2024  Through  10/27/2024"

# --- Weekly crime complaint table (rows 14-31) ---
# Row 14
$ws.Range("D14").Value = 1
$ws.Range("D14").NumberFormat = '#,##0'
$ws.Range("E14").Value = -100
$ws.Range("E14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J14").Value = 11
$ws.Range("K14").Value = -27.272727272727
$ws.Range("N14").Value = -81.395348837209

# Row 15
$ws.Range("C15").Value = 3
$ws.Range("F15").Value = 7
$ws.Range("I15").Value = 38
$ws.Range("K15").Value = 26.666666666666
$ws.Range("L15").Value = 58.333333333333
$ws.Range("M15").Value = 46.153846153846
$ws.Range("N15").Value = -54.761904761904

# Row 16
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = -12.5
$ws.Range("F16").Value = 21
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = 5
$ws.Range("I16").Value = 217
$ws.Range("J16").Value = 208
$ws.Range("K16").Value = 4.326923076923
$ws.Range("L16").Value = -19.926199261992
$ws.Range("M16").Value = -33.435582822085
$ws.Range("N16").Value = -89.084507042253

# Row 17
$ws.Range("C17").Value = 15
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = 50
$ws.Range("G17").Value = 51
$ws.Range("H17").Value = 27.450980392156
$ws.Range("I17").Value = 707
$ws.Range("J17").Value = 553
$ws.Range("K17").Value = 27.848101265822
$ws.Range("L17").Value = 30.202578268876
$ws.Range("M17").Value = 91.598915989159
$ws.Range("N17").Value = -32.858499525166

# Row 18
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 66.666666666666
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 162
$ws.Range("J18").Value = 158
$ws.Range("K18").Value = 2.531645569620
$ws.Range("L18").Value = -8.474576271186
$ws.Range("M18").Value = -54.494382022471
$ws.Range("N18").Value = -90.542907180385

# Row 19
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -46.153846153846
$ws.Range("F19").Value = 41
$ws.Range("G19").Value = 56
$ws.Range("H19").Value = -26.785714285714
$ws.Range("I19").Value = 475
$ws.Range("J19").Value = 551
$ws.Range("K19").Value = -13.793103448275
$ws.Range("L19").Value = -16.666666666666
$ws.Range("M19").Value = 3.711790393013
$ws.Range("N19").Value = -26.127527216174

# Row 20
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -60
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = 7.142857142857
$ws.Range("I20").Value = 166
$ws.Range("J20").Value = 180
$ws.Range("K20").Value = -7.777777777777
$ws.Range("L20").Value = -19.024390243902
$ws.Range("M20").Value = -29.059829059829
$ws.Range("N20").Value = -90.101371496720

# Row 21
$ws.Range("C21").Value = 39
$ws.Range("D21").Value = 40
$ws.Range("E21").Value = -2.5
$ws.Range("F21").Value = 161
$ws.Range("G21").Value = 154
$ws.Range("H21").Value = 4.545454545454
$ws.Range("I21").Value = 1773
$ws.Range("J21").Value = 1691
$ws.Range("K21").Value = 4.849201655824
$ws.Range("L21").Value = -1.827242524916
$ws.Range("M21").Value = -0.838926174496
$ws.Range("N21").Value = -75.378419663935

# Row 24
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = -53.571428571428
$ws.Range("F24").Value = 81
$ws.Range("G24").Value = 81
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 982
$ws.Range("J24").Value = 949
$ws.Range("K24").Value = 3.477344573234
$ws.Range("L24").Value = 9.476031215161
$ws.Range("M24").Value = 25.095541401273

# Row 25
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 19
$ws.Range("H25").Value = 72.727272727272
$ws.Range("I25").Value = 155
$ws.Range("J25").Value = 168
$ws.Range("K25").Value = -7.738095238095
$ws.Range("L25").Value = -23.267326732673

# Row 26
$ws.Range("C26").Value = 19
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = 35.714285714285
$ws.Range("F26").Value = 72
$ws.Range("H26").Value = 60
$ws.Range("I26").Value = 812
$ws.Range("J26").Value = 621
$ws.Range("K26").Value = 30.756843800322
$ws.Range("L26").Value = 30.756843800322
$ws.Range("M26").Value = 6.282722513089

# Row 27
$ws.Range("C27").Value = 3
$ws.Range("F27").Value = 7
$ws.Range("I27").Value = 52
$ws.Range("K27").Value = 30
$ws.Range("L27").Value = 33.333333333333

# Row 28
$ws.Range("D28").Value = 3
$ws.Range("D28").NumberFormat = '#,##0'
$ws.Range("E28").Value = -33.333333333333
$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F28").Value = 6
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 58
$ws.Range("J28").Value = 59
$ws.Range("K28").Value = -1.694915254237
$ws.Range("L28").Value = -9.375

# Row 29
$ws.Range("C29").Value = 1
$ws.Range("C29").NumberFormat = '#,##0'
$ws.Range("F29").Value = 1
$ws.Range("F29").NumberFormat = '#,##0'
$ws.Range("I29").Value = 27
$ws.Range("K29").Value = -20.588235294117
$ws.Range("L29").Value = -40
$ws.Range("M29").Value = -57.8125
$ws.Range("N29").Value = -86.363636363636

# Row 30
$ws.Range("C30").Value = 1
$ws.Range("C30").NumberFormat = '#,##0'
$ws.Range("F30").Value = 1
$ws.Range("F30").NumberFormat = '#,##0'
$ws.Range("I30").Value = 23
$ws.Range("K30").Value = -20.689655172413
$ws.Range("L30").Value = -39.473684210526
$ws.Range("M30").Value = -58.928571428571
$ws.Range("N30").Value = -87.150837988826

# Row 31
$ws.Range("C31").Value = 1
$ws.Range("C31").NumberFormat = '#,##0'
$ws.Range("F31").Value = 1
$ws.Range("F31").NumberFormat = '#,##0'
$ws.Range("I31").Value = 6
$ws.Range("K31").Value = 500
$ws.Range("L31").Value = 100
